# Rename the inline picture "name" (wp:docPr/@name, mirrored onto
# pic:cNvPr/@name) for the three logo images that live in the document's
# headers/footers:
#   - BTec logo (header, first-page):      image1.jpg -> image2.jpg
#   - Pearson logo (footer, default):      image2.png -> image1.png
#   - Pearson logo (footer, first-page):   image2.png -> image1.png
#
# InlineShape does not expose a settable "Name" property (matching real
# Word's object model), so each picture is temporarily converted to a
# floating Shape -- which *does* expose .Name -- renamed, then converted
# straight back to an inline picture so the wp:inline layout is preserved.

$d = $word.ActiveDocument

function Rename-InlinePicture($shape, $newName) {
    $floating = $shape.ConvertToShape()
    $floating.Name = $newName
    [void]$floating.ConvertToInlineShape()
}

for ($secIdx = 1; $secIdx -le $d.Sections.Count; $secIdx++) {
    $section = $d.Sections.Item($secIdx)

    for ($hfIdx = 1; $hfIdx -le 3; $hfIdx++) {
        $header = $section.Headers.Item($hfIdx)
        if ($header.Exists) {
            $shapes = $header.Range.InlineShapes
            for ($k = 1; $k -le $shapes.Count; $k++) {
                $shape = $shapes.Item($k)
                if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
                    Rename-InlinePicture $shape "image2.jpg"
                }
            }
        }

        $footer = $section.Footers.Item($hfIdx)
        if ($footer.Exists) {
            $shapes = $footer.Range.InlineShapes
            for ($k = 1; $k -le $shapes.Count; $k++) {
                $shape = $shapes.Item($k)
                if ($shape.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    Rename-InlinePicture $shape "image1.png"
                }
            }
        }
    }
}
